# Append new FMCG stock price history rows (2024-08-28 .. 2024-09-25) to Sheet1.
# The sheet's existing data runs through row 612 (date 2024-09-24); this appends
# rows 613-633 directly below it, extending the used range to A1:J633.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=613; A="2024-08-28"; C=661.9000244140625; D=1477.5; E=611.2000122070312; F=1308.050048828125; G=839.75; H=18182.55035400391; I=0; J=177.4530590978789},
    @{Row=614; A="2024-08-29"; C=660.75; D=1476.699951171875; E=603.6199951171875; F=1309.849975585938; G=813; H=18047.759765625; I=-0.007413183835854168; J=176.1375669485516},
    @{Row=615; A="2024-08-30"; C=647.1500244140625; D=1481.199951171875; E=600.3599853515625; F=1311.599975585938; G=813.2000122070312; H=17962.32995605469; I=-0.004733540931380745; J=175.3038125658468},
    @{Row=616; A="2024-09-02"; C=650.9500122070312; D=1470.050048828125; E=608.5800170898438; F=1316.800048828125; G=811.2000122070312; H=17987.74047851562; I=0.001414656256905703; J=175.5518072011525},
    @{Row=617; A="2024-09-03"; C=640.0499877929688; D=1460.75; E=599.9400024414062; F=1341.949951171875; G=822.3499755859375; H=17977.66967773438; I=-0.0005598702512568748; J=175.4535209667462},
    @{Row=618; A="2024-09-04"; C=645.5999755859375; D=1475.300048828125; E=609; F=1327.75; G=824.2999877929688; H=18052.54992675781; I=0.004165181047695957; J=176.1843166470284},
    @{Row=619; A="2024-09-05"; C=643.9000244140625; D=1457.699951171875; E=602.1799926757812; F=1254.800048828125; G=835.4000244140625; H=17792.94024658203; I=-0.01438077619112318; J=173.6506494209415},
    @{Row=620; A="2024-09-06"; C=665.25; D=1443.449951171875; E=597.2999877929688; F=1256.849975585938; G=832.7000122070312; H=17880.34979248047; I=0.004912597057432855; J=174.5037250903082},
    @{Row=621; A="2024-09-09"; C=676; D=1492.050048828125; E=610.3400268554688; F=1225.25; G=827.5999755859375; H=18025.32012939453; I=0.008107802061849448; J=175.9185667523958},
    @{Row=622; A="2024-09-10"; C=680; D=1503.050048828125; E=608; F=1246; G=824.75; H=18130.15014648438; I=0.005815709032478913; J=176.9416579500384},
    @{Row=623; A="2024-09-11"; C=680.4500122070312; D=1499.949951171875; E=627.6599731445312; F=1229; G=814; H=18088.97985839844; I=-0.002270818926114677; J=176.5398554843474},
    @{Row=624; A="2024-09-12"; C=686.0999755859375; D=1513.449951171875; E=645.5999755859375; F=1224.849975585938; G=809.7000122070312; H=18193.19958496094; I=0.00576150382046627; J=177.556990536185},
    @{Row=625; A="2024-09-13"; C=681.9500122070312; D=1491.300048828125; E=646.6500244140625; F=1229.300048828125; G=788.0499877929688; H=18027.60040283203; I=-0.009102257211853799; J=175.9408211385619},
    @{Row=626; A="2024-09-16"; C=695.2000122070312; D=1456.349975585938; E=621.0499877929688; F=1219.699951171875; G=751.9500122070312; H=17765.49987792969; I=-0.01453884704817227; J=173.3828444504985},
    @{Row=627; A="2024-09-17"; C=692; D=1459.400024414062; E=649.6500244140625; F=1222.949951171875; G=746.75; H=17827; I=0.003461772676980224; J=173.9830564440744},
    @{Row=628; A="2024-09-18"; C=695.2999877929688; D=1432.150024414062; E=646.7000122070312; F=1224.550048828125; G=744.5999755859375; H=17755.70007324219; I=-0.003999547133999692; J=173.287203009309},
    @{Row=629; A="2024-09-19"; C=697; D=1444.849975585938; E=652.1500244140625; F=1197.849975585938; G=747.2000122070312; H=17752.34997558594; I=-0.0001886773060161447; J=173.2545076466781},
    @{Row=630; A="2024-09-20"; C=709; D=1456.599975585938; E=654.4500122070312; F=1206.300048828125; G=747.5499877929688; H=17905.25006103516; I=0.008612949027001824; J=174.7467398897373},
    @{Row=631; A="2024-09-23"; C=702.5; D=1449.300048828125; E=654.0999755859375; F=1190; G=763.75; H=17852.70007324219; I=-0.00293489270542646; J=174.2338769575378},
    @{Row=632; A="2024-09-24"; C=705.0999755859375; D=1446.349975585938; E=646.8499755859375; F=1194.699951171875; G=760.9500122070312; H=17843.19958496094; I=-0.0005321597429113499; J=174.1411567023696},
    @{Row=633; A="2024-09-25"; C=689.2000122070312; D=1429.550048828125; E=633.2999877929688; F=1175.349975585938; G=742.5499877929688; H=17509.20007324219; I=-0.0187185885652627; J=170.8814800377791}
)

$firstNewRow = 613
$lastNewRow = 633

# Column I (9) carries a column-wide percentage number format ("0.0%") defined on
# the worksheet's <col> entry for that column. A brand-new cell written into that
# column with a plain value assignment would inherit this column default, but the
# source data for these rows is unformatted (same as the existing I593:I612
# cells). Pre-seed the new I-column cells' formatting from the row immediately
# above (which already has no special format) so the appended cells stay plain.
$ws.Cells.Item($firstNewRow - 1, 9).Copy()
$ws.Range($ws.Cells.Item($firstNewRow, 9), $ws.Cells.Item($lastNewRow, 9)).PasteSpecial(-4122)
$excel.CutCopyMode = $false

foreach ($d in $data) {
    $row = $d.Row

    # Column A holds a "YYYY-MM-DD" date string stored as literal text (matching
    # every other row in the sheet); force text formatting first so Excel does
    # not auto-convert the assignment into a date serial number.
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $d.A

    $ws.Cells.Item($row, 3).Value = $d.C
    $ws.Cells.Item($row, 4).Value = $d.D
    $ws.Cells.Item($row, 5).Value = $d.E
    $ws.Cells.Item($row, 6).Value = $d.F
    $ws.Cells.Item($row, 7).Value = $d.G
    $ws.Cells.Item($row, 8).Value = $d.H
    $ws.Cells.Item($row, 9).Value = $d.I
    $ws.Cells.Item($row, 10).Value = $d.J
}
